# Auto-generated Excel COM-interop edit script
# Refreshes the crypto table's Price (D) and Volume(1h) (E) columns
# per the scraped source update (GitHub Actions cron job).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D/E hold numeric-looking text (e.g. '304.40', '46.918.82',
# '  +6.12%  ') that must stay plain text, matching the original
# inline-string cells. Force text format first so Excel's COM layer
# doesn't auto-convert values like '304.40' into the number 304.4,
# then clear the temporary formatting so no stray style is left behind.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '46.888.09'
$ws.Range('E2').Value = '  +6.27%  '
$ws.Range('D3').Value = '2.329.96'
$ws.Range('E3').Value = '  +5.07%  '
$ws.Range('D5').Value = '304.40'
$ws.Range('E5').Value = '  +1.06%  '
$ws.Range('D6').Value = '97.49'
$ws.Range('E6').Value = '  +9.94%  '
$ws.Range('D7').Value = '0.577'
$ws.Range('E7').Value = '  +4.26%  '
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.56%  '
$ws.Range('E9').Value = '  +9.82%  '
$ws.Range('D10').Value = '36.05'
$ws.Range('E10').Value = '  +7.82%  '
$ws.Range('D11').Value = '0.0812'
$ws.Range('E11').Value = '  +4.17%  '
$ws.Range('E12').Value = '  +8.86%  '
$ws.Range('E13').Value = '  +0.93%  '
$ws.Range('D14').Value = '2.683.57'
$ws.Range('E14').Value = '  +5.03%  '
$ws.Range('D15').Value = '2.329.41'
$ws.Range('E15').Value = '  +1.79%  '
$ws.Range('D16').Value = '0.842'
$ws.Range('E16').Value = '  +5.50%  '
$ws.Range('D17').Value = '14.13'
$ws.Range('E17').Value = '  +8.03%  '
$ws.Range('D18').Value = '46.770.23'
$ws.Range('E18').Value = '  +6.48%  '
$ws.Range('D19').Value = '13.70'
$ws.Range('E19').Value = '  +22.44%  '
$ws.Range('D20').Value = '0.0₃0956'
$ws.Range('E21').Value = '  +4.24%  '
$ws.Range('D22').Value = '67.82'
$ws.Range('E22').Value = '  +5.85%  '
$ws.Range('D23').Value = '254.46'
$ws.Range('E23').Value = '  +9.46%  '
$ws.Range('E24').Value = '  +4.63%  '
$ws.Range('E25').Value = '  +6.32%  '
$ws.Range('E26').Value = '  -0.08%  '
$ws.Range('D27').Value = '42.11'
$ws.Range('E27').Value = '  +16.88%  '
$ws.Range('D28').Value = '2.31'
$ws.Range('E28').Value = '  +2.69%  '
$ws.Range('D29').Value = '9.96'
$ws.Range('E29').Value = '  +6.27%  '
$ws.Range('D30').Value = '20.30'
$ws.Range('E30').Value = '  +4.97%  '
$ws.Range('E31').Value = '  +5.33%  '
$ws.Range('D32').Value = '0.0817'
$ws.Range('E32').Value = '  +8.62%  '
$ws.Range('D33').Value = '147.99'
$ws.Range('E33').Value = '  +1.55%  '
$ws.Range('D34').Value = '2.65'
$ws.Range('E34').Value = '  +2.32%  '
$ws.Range('E35').Value = '  +7.81%  '
$ws.Range('D36').Value = '0.113'
$ws.Range('E36').Value = '  +7.60%  '
$ws.Range('E37').Value = '  +3.69%  '
$ws.Range('E38').Value = '  +6.53%  '
$ws.Range('D39').Value = '4.06'
$ws.Range('E39').Value = '  +12.98%  '
$ws.Range('E40').Value = '  +9.12%  '
$ws.Range('D41').Value = '3.40'
$ws.Range('E41').Value = '  +6.91%  '
$ws.Range('D42').Value = '14.13'
$ws.Range('E42').Value = '  -0.53%  '
$ws.Range('E43').Value = '  -0.97%  '
$ws.Range('E44').Value = '  +20.05%  '
$ws.Range('D45').Value = '92.23'
$ws.Range('D46').Value = '1.805.26'
$ws.Range('E46').Value = '  +4.12%  '
$ws.Range('E47').Value = '  +9.34%  '
$ws.Range('D48').Value = '74.07'
$ws.Range('E48').Value = '  +11.93%  '
$ws.Range('D49').Value = '98.96'
$ws.Range('E49').Value = '  +4.75%  '
$ws.Range('D50').Value = '4.90'
$ws.Range('E50').Value = '  +7.07%  '
$ws.Range('D51').Value = '55.29'
$ws.Range('E51').Value = '  +7.21%  '

# Drop the temporary text-number-format so cells keep their original
# (default/general) style index, matching the source workbook.
$dataRange.ClearFormats()
